$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TASKS")

# Select cell I10 on the TASKS sheet (mirrors the sheetView selection change)
$ws.Activate()
$ws.Range("I10").Select()

# E22: update the date value (45615 -> 45617, i.e. 19/11/2024 -> 21/11/2024)
$ws.Range("E22").Value = 45617

# Row 27: responsible changed to "Walid", state changed to "EC", date filled in
$ws.Range("C27").Value = "Walid"
$ws.Range("D27").Value = "EC"
$ws.Range("E27").Value = 45617

# Rows 29 & 30: state changed from "EC" to "V"
$ws.Range("D29").Value = "V"
$ws.Range("D30").Value = "V"
